$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two new SICK LEAVE (SL(1-0-0)) entries were recorded, requiring two new rows
# to be inserted into the leave-card table right before the (previously) next
# blank period row (old row 88, dated 8/1/2023).

# 1) Insert two blank rows at the insertion point (row 88).
$ws.Rows("88:89").Insert()

# 2) Seed the two new rows with the formatting/formula pattern used by the
#    other "SL(1-0-0)" leave rows already in the table (row 82 is one).
$ws.Range("A82:K82").Copy($ws.Range("A88:K88"))
$ws.Range("A82:K82").Copy($ws.Range("A89:K89"))

# 3) Re-assert the calculated "EARNED " column formula (table structured
#    reference) on the two new rows.
$ws.Range("G88").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G89").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# 4) Fill in the actual data for the two new rows: 1 day "Absence Undertime
#    W/ Pay" each, dated in the REMARKS column.
$ws.Range("H88").Value = 1
$ws.Range("K88").Value = 45124

$ws.Range("H89").Value = 1
$ws.Range("K89").Value = 45120

# 5) Grow Table1 so it covers the two freshly-inserted rows plus the two
#    rows the insert pushed past the table's old bottom edge (A8:K135 -> A8:K137).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A8:K137"))

# 6) The table resize re-creates the calculated column formula on the rows
#    that now trail the old table end; make sure they keep the same
#    structured-reference form as the rest of the column.
$ws.Range("G136").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G137").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

$excel.Calculate()

# 7) Leave the selection where the editor's cursor ended up.
$ws.Range("I89").Select()
